$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.311.12"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "'1.879.29"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").Value = "'246.24"
$ws.Range("E5").Value = "  -3.07%  "
$ws.Range("D6").Value = "'0.684"
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("D8").Value = "'43.58"
$ws.Range("E8").Value = "  +4.82%  "
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").Value = "'53.74"
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("E11").Value = "  -2.87%  "
$ws.Range("D12").Value = "'0.0973"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").Value = "'2.152.26"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "'0.763"
$ws.Range("E15").Value = "  +3.81%  "
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("D17").Value = "'1.889.76"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "'35.333.39"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").Value = "'72.75"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("D20").Value = "'0.0₃0821"
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("D21").Value = "'243.98"
$ws.Range("D22").Value = "'12.85"
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("D23").Value = "'4.96"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("D24").Value = "'2.63"
$ws.Range("E24").Value = "  +7.10%  "
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("E26").Value = "  -6.00%  "
$ws.Range("D27").Value = "'165.82"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "'8.53"
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").Value = "'4.128.46"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +5.73%  "
$ws.Range("D33").Value = "'2.04"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("D34").Value = "'4.30"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").Value = "'0.0593"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("E36").Value = "  -2.65%  "
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D40").Value = "'0.0724"
$ws.Range("E40").Value = "  +11.03%  "
$ws.Range("D41").Value = "'17.75"
$ws.Range("E41").Value = "  +4.12%  "
$ws.Range("D42").Value = "'0.0218"
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("D43").Value = "'96.03"
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").Value = "'1.303.03"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("E47").Value = "  +5.82%  "
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").Value = "'2.73"
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("D50").Value = "'11.93"
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("D51").Value = "'6.22"
$ws.Range("E51").Value = "  -5.69%  "
